# "Small changes: Do not print innecesary stuff"
# - Insert 2 new rows (HDC('normal', 15) and HDC('normal', 20)) into the
#   first correlation block, pushing the later blocks down by 2 rows.
# - Add formulas (AVERAGE of the underlying runs) to several existing cells
#   that used to hold a single literal number.
# - Highlight the "Mi dataset (64)" comparison block (rows 6-12) with a
#   grey fill, and mark a few standout numbers in bold purple.
# - Bold (no fill) the row-max values picked out in the HDC('normal', *) block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room for the two new rows (HDC('normal', 15) / HDC('normal', 20))
#    by inserting before the old row 17 ("FALTARÁ 20?" merged banner).
# ---------------------------------------------------------------------
$ws.Rows("17:18").Insert()

# ---------------------------------------------------------------------
# 2) Update the formulas/values for the HDC('normal', n) block (rows 13-18)
# ---------------------------------------------------------------------
$ws.Range("F13").Formula = "=AVERAGE(0.6091,0.5923)"
$ws.Range("G13").Formula = "=AVERAGE(0.5912)"
$ws.Range("H13").Formula = "=AVERAGE(0.513098,0.5024881,0.5037)"

$ws.Range("F14").Formula = "=AVERAGE(0.6204,0.644)"
$ws.Range("G14").Formula = "=AVERAGE(0.639)"
$ws.Range("H14").Formula = "=AVERAGE(0.4969,0.4987)"

$ws.Range("G15").Formula = "=AVERAGE(0.65324)"
$ws.Range("H15").Formula = "=AVERAGE(0.5332,0.5203)"

$ws.Range("F16").Formula = "=AVERAGE(0.73065)"
$ws.Range("G16").Formula = "=AVERAGE(0.6565)"
$ws.Range("H16").Formula = "=AVERAGE(0.5793,0.5837)"

# New row 17: HDC('normal', 15)
$ws.Range("E17").Value2 = "HDC('normal', 15)"
$ws.Range("F17").Formula = "=AVERAGE(0.7227)"
$ws.Range("G17").Formula = "=AVERAGE(0.6609)"
$ws.Range("H17").Formula = "=AVERAGE(0.6462)"

# New row 18: HDC('normal', 20)
$ws.Range("E18").Value2 = "HDC('normal', 20)"
$ws.Range("F18").Formula = "=AVERAGE(0.70125)"
$ws.Range("G18").Formula = "=AVERAGE(0.6568)"
$ws.Range("H18").Formula = "=AVERAGE(0.6554)"

# New rows inherited row-16's formatting via the insert; re-apply the
# label style explicitly so E17/E18 match the other labels in the block.
$ws.Range("E17:E18").HorizontalAlignment = -4108
$ws.Range("E17:E18").VerticalAlignment = -4130
$ws.Range("E17:E18").Font.Bold = $true
$ws.Range("F17:H18").HorizontalAlignment = -4108
$ws.Range("F17:H18").VerticalAlignment = -4108

# Row-max picks in the HDC('normal', n) block get bold (no fill), matching
# the pre-existing "winner" styling already used elsewhere in the sheet.
$ws.Range("F16").Font.Bold = $true
$ws.Range("F16").VerticalAlignment = -4108
$ws.Range("F17").Font.Bold = $true
$ws.Range("F17").VerticalAlignment = -4108
$ws.Range("H18").Font.Bold = $true
$ws.Range("H18").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# 3) Highlight the top comparison block (rows 6-12) with a grey fill.
# ---------------------------------------------------------------------
# Header row (column titles)
$hdr = $ws.Range("F6:H6")
$hdr.Interior.Color = 14277081
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108

# Row labels (E7:E12)
$lbl = $ws.Range("E7:E12")
$lbl.Interior.Color = 14277081
$lbl.Font.Bold = $true
$lbl.HorizontalAlignment = -4108

# Data block (F7:H12)
$data = $ws.Range("F7:H12")
$data.Interior.ThemeColor = 4
$data.HorizontalAlignment = -4108
$data.VerticalAlignment = -4108

# Stand-out values: bold purple font, keep the grey fill.
$purple = 10498160
foreach ($addr in @("F7", "H7", "G8")) {
    $c = $ws.Range($addr)
    $c.Font.Bold = $true
    $c.Font.Color = $purple
    $c.Interior.ThemeColor = 4
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108
}

# ---------------------------------------------------------------------
# 4) View bookkeeping to mirror what Excel records after this kind of
#    edit (selection sitting on the freshly edited merged banner row).
# ---------------------------------------------------------------------
$ws.Range("E19:H19").Select()
$excel.ActiveWindow.ScrollRow = 4
